$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New short entry far below the table (row 19): just a YouTube link.
# Entered first so the shared-string table picks up its text before the
# longer log-entry row below (matches the order new strings were appended
# to the workbook when it was authored).
$ws.Range("B19").Value = "https://www.youtube.com/watch?v=z0MimkXIvE8"

# New log row 15, same look & feel as the other "grey" rows (e.g. row 13):
# copy the formatting from row 13 first, then fill in the new content.
$ws.Range("B13:D13").Copy()
$ws.Range("B15:D15").PasteSpecial(-4122)

$ws.Range("B15").Value = "Extra videos over postman bekeken en een paar extra requests gemaakt om mijn spotify playlist bij te vullen"
$ws.Range("C15").Value = "11/15/2021"
$ws.Range("D15").Value = "65 minuten"

# The extra wrapped text needs a taller row.
$ws.Rows.Item(15).RowHeight = 60

# Leave the selection where the user last clicked.
[void]$ws.Range("D16").Select()
